# Updated symbol list on Tue Dec 13 14:16:10 UTC 2022 with GitHub Actions
#
# For every coin row (2..51) on Sheet1:
#   - column G ("Hora") flips from "13" to "14"
#   - column D ("Price") is refreshed to the new quote, where available
#     (rows whose price is currently unavailable, e.g. "--", keep their
#     existing value and only the hour ticks forward)
#
# Values are written as literal text (not numbers) so that things like the
# trailing zero in "0.06290" or "6.370" survive exactly as scraped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new Price (column D); rows absent from this map had no new price
# pulled (the source cell already reads "--") and are left untouched.
$newPrices = @{
    2  = "274.07"
    3  = "21.57"
    4  = "6.370"
    5  = "0.06290"
    6  = "3.616"
    7  = "6.608"
    8  = "1.397"
    9  = "0.8287"
    10 = "0.01386"
    11 = "0.1589"
    12 = "0.08371"
    13 = "0.03458"
    14 = "0.03212"
    15 = "4.083"
    16 = "0.09281"
    17 = "0.001651"
    18 = "0.04742"
    19 = "0.006298"
    20 = "0.005993"
    21 = "0.001073"
    22 = "0.0001495"
    23 = "3.726"
    24 = "2.415"
    25 = "0.3329"
    26 = "0.1258"
    28 = "0.0002694"
    40 = "0.04745"
    41 = "0.007066"
    42 = "0.1172"
    43 = "0.003369"
    44 = "0.01178"
    45 = "0.00006073"
    46 = "0.0009864"
    47 = "0.00000000747"
    48 = "0.7793"
    49 = "0.002461"
    50 = "0.00001295"
    51 = "0.01236"
}

function Set-TextValue($cell, [string]$value) {
    # Force text storage so numeric-looking strings ("14", "0.06290", ...)
    # keep their exact textual form instead of being normalised into a
    # floating point number, then drop the formatting footprint again so
    # the cell's style stays identical to how it started.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

for ($row = 2; $row -le 51; $row++) {
    Set-TextValue $ws.Cells.Item($row, 7) "14"

    if ($newPrices.ContainsKey($row)) {
        Set-TextValue $ws.Cells.Item($row, 4) $newPrices[$row]
    }
}
